$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.791.32"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "2.285.57"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.82"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.645"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.08"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +6.77%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.645"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.60"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.50"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "2.629.78"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.12"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.872"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "2.288.89"
$ws.Range("E17").Value = "  +1.66%  "
$ws.Range("D18").Value = "42.695.33"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.41"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.28"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("E23").Value = "  +5.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.86"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.34"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.40"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.17"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.56"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.04"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0876"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +11.78%  "
$ws.Range("E32").Value = "  -3.27%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.23"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("E35").Value = "  +2.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.60"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.79"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0304"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.70"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +9.80%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +2.10%  "
$ws.Range("E42").Value = "  +5.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.10"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.00"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.83"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.30"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +11.30%  "
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.17"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.23"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.08%  "
